# issue #5: stock data output to json file
#
# Adds a "property_category" column (with value "stock") to the 股票
# (stock) worksheet, between the existing "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票 (stock) sheet

# Insert a new column at H, pushing date/legislator_name/legislator_id
# (and the K2 amount) one column to the right.
$ws.Columns("H").Insert()

# Populate the new column's header and value.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
